$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '33.039.97'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +10.56%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.756.74'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +5.93%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.991'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.70%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.79'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +4.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.544'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +4.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.991'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.65%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.71'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +9.37%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.42'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.76%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.276'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +5.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0661'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +7.99%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0915'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.004.64'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +5.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.755.77'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +5.90%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '10.62'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +6.05%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.635'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +5.78%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.27'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +8.37%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '33.015.48'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +10.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '68.50'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +5.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '259.24'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +7.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0737'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +4.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.991'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.54%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.43'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +4.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.34'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.17'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.43'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.51'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +4.99%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.95'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.992'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.86'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +13.75%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0517'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.18'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +5.68%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +7.40%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.554.15'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +8.53%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.79'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +4.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '85.44'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +9.36%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.67%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.628'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +9.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0185'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +5.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.72'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.30'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.25%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +7.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.869'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0509'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.37%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '55.76'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +10.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.06'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +5.79%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.898.53'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +5.33%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.66'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +5.56%  '
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.990'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.79%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.12'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +21.74%  '
